$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.001.44'
$ws.Range('E2').Value = '  +7.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.249.44'
$ws.Range('E3').Value = '  +3.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '394.34'
$ws.Range('E5').Value = '  -0.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.57'
$ws.Range('E6').Value = '  +2.22%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.247.42'
$ws.Range('E7').Value = '  +2.95%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  +3.78%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.616'
$ws.Range('E10').Value = '  +1.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '38.86'
$ws.Range('E11').Value = '  +0.98%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0982'
$ws.Range('E12').Value = '  +13.03%  '

$ws.Range('E13').Value = '  +1.95%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.758.59'
$ws.Range('E14').Value = '  +2.89%  '

$ws.Range('E15').Value = '  +1.42%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.95'
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.237.28'
$ws.Range('E17').Value = '  +2.89%  '

$ws.Range('E18').Value = '  -2.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  -0.39%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '56.790.12'
$ws.Range('E20').Value = '  +6.84%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.32'
$ws.Range('E21').Value = '  +1.91%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000106'
$ws.Range('E22').Value = '  +8.37%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.00'
$ws.Range('E23').Value = '  +1.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '297.44'
$ws.Range('E24').Value = '  +9.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.57'
$ws.Range('E25').Value = '  +3.43%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.14'
$ws.Range('E26').Value = '  -2.28%  '

$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.39'
$ws.Range('E27').Value = '  +3.34%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '27.92'
$ws.Range('E28').Value = '  +0.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.73'
$ws.Range('E29').Value = '  -4.45%  '

$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.168'
$ws.Range('E30').Value = '  -1.41%  '

$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.20'
$ws.Range('E31').Value = '  -3.95%  '

$ws.Range('E32').Value = '  +0.00%  '

$ws.Range('E33').Value = '  -0.36%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.96'
$ws.Range('E34').Value = '  -0.65%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '37.14'
$ws.Range('E35').Value = '  -0.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0483'
$ws.Range('E36').Value = '  -1.68%  '

$ws.Range('E37').Value = '  +1.49%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.55'
$ws.Range('E38').Value = '  +2.16%  '

$ws.Range('E39').Value = '  +0.68%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.24%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.06'
$ws.Range('E41').Value = '  +11.35%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '134.11'
$ws.Range('E42').Value = '  +3.20%  '

$ws.Range('E43').Value = '  +0.86%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.119'
$ws.Range('E44').Value = '  +2.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.99'
$ws.Range('E45').Value = '  -1.49%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.93'
$ws.Range('E46').Value = '  -6.12%  '

$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.283'
$ws.Range('E47').Value = '  -3.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.86'
$ws.Range('E48').Value = '  -1.91%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.143.33'
$ws.Range('E49').Value = '  +2.69%  '

$ws.Range('E50').Value = '  -0.55%  '

$ws.Range('E51').Value = '  +23.30%  '
